$p = $ppt.ActivePresentation
$nl = [char]13

# Slide 2 ("Объект 3" placeholder, shape #2): the first bullet was split
# across two runs ("Описани" + "е игры"); retype it as a single run
# "Описание игры", keeping the other bullets untouched.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
# The combined plain text is identical to what's already there ("Описани" +
# "е игры" == "Описание игры"), so a direct re-assignment is a no-op for the
# engine's change-detection. Flip through a distinct placeholder string first
# to force the run structure to actually collapse to one run, then set the
# final text.
$tr2.Text = "~PLACEHOLDER~" + $nl + "Использованные библиотеки" + $nl + "Описание кода" + $nl + "Вывод"
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Описание игры" + $nl + "Использованные библиотеки" + $nl + "Описание кода" + $nl + "Вывод"

# Slide 3 title ("Заголовок 1", shape #1): "Наши обязательства" -> "Сюжет"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Сюжет"
